$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Rows.Item(75).Insert()
$ws.Range("A75").Value = "Hello Title"
$ws.Range("A75").Font.Bold = $true
$ws.Range("A75").Font.Size = 16
$ws.Range("A75:N75").Interior.ThemeColor = 5
